$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.000.60"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "3.924.66"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'483.63"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "'146.65"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  +7.85%  "
$ws.Range("E11").Value = "  +13.28%  "
$ws.Range("D12").Value = "'42.61"
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.44"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.543.92"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "3.940.30"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'14.57"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'19.67"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "69.082.41"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "'434.29"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'14.59"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'3.35"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").Value = "'87.95"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  +15.37%  "
$ws.Range("D26").Value = "'3.56"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "'38.32"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "'10.35"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "'5.92"
$ws.Range("E29").Value = "  +8.40%  "
$ws.Range("D30").Value = "'709.09"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.129"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'13.24"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "'2.85"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "0.0₃0941"
$ws.Range("E34").Value = "  +35.59%  "
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").Value = "'58.70"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("E37").Value = "  -7.02%  "
$ws.Range("D38").Value = "'5.58"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'0.0471"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  +8.91%  "
$ws.Range("D42").Value = "'3.01"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("D43").Value = "'2.99"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").Value = "'0.341"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "'148.07"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "'3.12"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").Value = "'2.83"
$ws.Range("E51").Value = "  -2.14%  "
